$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '30.699.32'
$ws.Cells.Item(2, 5).Value = '  +1.19%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.880.10'
$ws.Cells.Item(3, 5).Value = '  -0.03%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9994'
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '239.05'
$ws.Cells.Item(5, 5).Value = '  +0.53%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.9995'
$ws.Cells.Item(6, 5).Value = '  +0.00%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4809'
$ws.Cells.Item(7, 5).Value = '  -0.39%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2841'
$ws.Cells.Item(8, 5).Value = '  -1.89%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06541'
$ws.Cells.Item(9, 5).Value = '  -0.75%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.896.77'
$ws.Cells.Item(10, 5).Value = '  +0.92%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07482'
$ws.Cells.Item(11, 5).Value = '  +1.27%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '16.64'
$ws.Cells.Item(12, 5).Value = '  -1.76%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.125'
$ws.Cells.Item(13, 5).Value = '  -0.98%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '88.87'
$ws.Cells.Item(14, 5).Value = '  +0.72%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.6676'
$ws.Cells.Item(15, 5).Value = '  +1.01%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '30.651.31'
$ws.Cells.Item(16, 5).Value = '  +1.16%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '13.39'
$ws.Cells.Item(17, 5).Value = '  -1.01%  '
$ws.Cells.Item(18, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.254.79'
$ws.Cells.Item(18, 5).Value = '  +5.64%  '
$ws.Cells.Item(19, 2).Value = 'Dai'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '1.000'
$ws.Cells.Item(19, 5).Value = '  +0.08%  '
$ws.Cells.Item(20, 2).Value = 'ShibaInu'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.000007646'
$ws.Cells.Item(20, 5).Value = '  -1.34%  '
$ws.Cells.Item(21, 2).Value = 'BitcoinCash'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '234.33'
$ws.Cells.Item(21, 5).Value = '  +17.35%  '
$ws.Cells.Item(22, 5).Value = '  -2.21%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.9995'
$ws.Cells.Item(23, 5).Value = '  -0.01%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.235'
$ws.Cells.Item(24, 5).Value = '  +1.07%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.337'
$ws.Cells.Item(25, 5).Value = '  -0.94%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '167.14'
$ws.Cells.Item(26, 5).Value = '  +1.73%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '18.73'
$ws.Cells.Item(27, 5).Value = '  +2.60%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.964'
$ws.Cells.Item(28, 5).Value = '  +1.82%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.455'
$ws.Cells.Item(29, 5).Value = '  +1.37%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.09585'
$ws.Cells.Item(30, 5).Value = '  +4.90%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.366'
$ws.Cells.Item(31, 5).Value = '  +2.44%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.048'
$ws.Cells.Item(32, 5).Value = '  +0.10%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.05050'
$ws.Cells.Item(33, 5).Value = '  +0.10%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.214'
$ws.Cells.Item(34, 5).Value = '  +5.76%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.7520'
$ws.Cells.Item(35, 5).Value = '  +1.33%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.706'
$ws.Cells.Item(36, 5).Value = '  +0.01%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.01862'
$ws.Cells.Item(37, 5).Value = '  +1.13%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.631'
$ws.Cells.Item(38, 5).Value = '  -0.04%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.9185'
$ws.Cells.Item(39, 5).Value = '  +0.50%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.086'
$ws.Cells.Item(40, 5).Value = '  +0.46%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '106.09'
$ws.Cells.Item(41, 5).Value = '  -0.39%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.4294'
$ws.Cells.Item(42, 5).Value = '  -0.63%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '5.820'
$ws.Cells.Item(43, 5).Value = '  -1.33%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.9996'
$ws.Cells.Item(44, 5).Value = '  +0.02%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '7.480'
$ws.Cells.Item(45, 5).Value = '  -2.25%  '
$ws.Cells.Item(46, 2).Value = 'Algorand'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.1292'
$ws.Cells.Item(46, 5).Value = '  -4.29%  '
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '64.59'
$ws.Cells.Item(47, 5).Value = '  -0.55%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.487'
$ws.Cells.Item(48, 5).Value = '  -4.63%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '8.993'
$ws.Cells.Item(49, 5).Value = '  +1.10%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '33.97'
$ws.Cells.Item(50, 5).Value = '  -0.49%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.3898'
$ws.Cells.Item(51, 5).Value = '  +0.51%  '
